$d = $word.ActiveDocument

# Change 1: Replace the domain walls paragraph text
$d.Content.Find.Execute(
    "We have rewritten the section on domain walls to make it clearer. We have also emphasized that the domain wall picture only provides a heuristic understanding of the critical dimension, and that the generalized Mermin-Wagner argument comes from correlation function calculations.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The domain wall discussion was meant to be a heuristic argument that would help with intuition prior to the correlation function calculations. Since it doesn’t seem to help, we have removed it.",
    2
)

# Change 2: Replace the IV.A/IV.B paragraph text and remove the red "DO WE WANT..." run
$d.Content.Find.Execute(
    "We have updated sections IV.A and IV.B to provide a clearer picture of the results of symmetry breaking in our models. In particular, section IV.B now says in which dimensions it is possible to break any maximal multipole group to any of its maximal subgroups. DO WE WANT TO COUNT GOLDSTONE MODES?",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "We have updated sections IV.A and IV.B to provide a clearer picture of the results of symmetry breaking in our models. In particular, section IV.B now says in which dimensions it is possible to break any maximal multipole group to any of its maximal subgroups. We have included a note on the subtleties of counting massless modes and of deciding which such modes are true Goldstone modes.",
    2
)

Write-Host "done"
